# Update course data on the "courses" sheet:
#  - Add a cricosCode value in column B for every data row (2-18)
#  - Replace the placeholder department value ("SHELDON SCHOOL OF HOSPITALITY")
#    in column C with the real department for each row
#  - Fix a spelling typo in the combined package name in D11
#    ("Commerical" -> "Commercial" for the first course in the list)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cricosCode (column B) per row
$cricosCodes = @{
    2  = "093480J"
    3  = "093481G"
    4  = "095244E"
    5  = "093479B"
    6  = "106672A"
    7  = "091074F"
    8  = "091132A"
    9  = "095245D"
    10 = "095246C"
    11 = "093480J/095244E/091074F"
    12 = "095244E/091132A"
    13 = "093481G/093479B/091074F"
    14 = "106672A/091132A"
    15 = "093479B/091132A"
    16 = "090975K/091074F"
    17 = "091074F/091132A"
    18 = "095245D/095246C"
}

# department (column C) per row
$departments = @{
    2  = "Cookery"
    3  = "Patisserie and Baking"
    4  = "Cookery"
    5  = "Patisserie and Baking"
    6  = "Patisserie and Baking"
    7  = "Hospitality"
    8  = "Hospitality"
    9  = "Travel and Tourism"
    10 = "Travel and Tourism"
    11 = "Packages"
    12 = "Packages"
    13 = "Packages"
    14 = "Packages"
    15 = "Packages"
    16 = "Packages"
    17 = "Packages"
    18 = "Packages"
}

# Rows 11-18 list package/combined courses and already use a wrapped-text
# style on column A; mirror that same wrap formatting onto the newly
# populated column B for those rows.
$wrapRows = 11..18

foreach ($r in 2..18) {
    $ws.Cells.Item($r, 2).Value = $cricosCodes[$r]
    if ($wrapRows -contains $r) {
        $ws.Cells.Item($r, 2).WrapText = $true
    }

    $ws.Cells.Item($r, 3).Value = $departments[$r]
}

# Correct "Commerical" -> "Commercial" in the first course name of the
# Certificate III/IV + Diploma package (row 11, column D).
$ws.Range("D11").Value = "Certificate III in Commercial Cookery + Certificate IV in Commerical Cookery + Diploma of Hospitality Management"
